$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.994.72'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.558.93'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '208.02'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.03'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0596'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0855'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '1.783.59'
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = '1.560.96'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.72'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.519'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '26.996.59'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '0.0₃0706'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '215.72'
$ws.Range('D19').Style = "Normal"
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.13'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.93'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.84'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.59'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.105'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0473'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('E31').Value = '  +3.58%  '
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  +3.56%  '
$ws.Range('D34').Value = '1.425.52'
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +9.30%  '
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('E39').Value = '  +2.21%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.88'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.807'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.54'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').Value = '1.696.87'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.83'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').Value = '0.0₆0104'
$ws.Range('E49').Value = '  +3.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0517'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0956'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.53%  '
